$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.028.17'
$ws.Range('E2').Value = '  -0.43%  '
$ws.Range('D3').Value = '1.829.55'
$ws.Range('E3').Value = '  -0.13%  '
$ws.Range('D4').Value = '0.9989'
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').Value = '241.08'
$ws.Range('E5').Value = '  -0.28%  '
$ws.Range('D6').Value = '0.6240'
$ws.Range('E6').Value = '  -5.11%  '
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('B8').Value = 'OKB'
$ws.Range('C8').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D8').Value = '44.72'
$ws.Range('E8').Value = '  +6.98%  '
$ws.Range('B9').Value = 'Dogecoin'
$ws.Range('C9').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D9').Value = '0.07528'
$ws.Range('E9').Value = '  +1.84%  '
$ws.Range('D10').Value = '0.2910'
$ws.Range('E10').Value = '  -0.29%  '
$ws.Range('D11').Value = '22.77'
$ws.Range('E12').Value = '  -1.22%  '
$ws.Range('D13').Value = '1.829.22'
$ws.Range('E13').Value = '  -0.70%  '
$ws.Range('D14').Value = '4.955'
$ws.Range('E14').Value = '  -0.63%  '
$ws.Range('D15').Value = '0.6637'
$ws.Range('E15').Value = '  -0.17%  '
$ws.Range('D16').Value = '82.30'
$ws.Range('E16').Value = '  -0.50%  '
$ws.Range('D17').Value = '0.000009069'
$ws.Range('E17').Value = '  +7.65%  '
$ws.Range('D18').Value = '6.005'
$ws.Range('E18').Value = '  -1.91%  '
$ws.Range('D19').Value = '28.938.75'
$ws.Range('E19').Value = '  -0.77%  '
$ws.Range('B20').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C20').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D20').Value = '2.082.66'
$ws.Range('E20').Value = '  -0.50%  '
$ws.Range('B21').Value = 'BitcoinCash'
$ws.Range('C21').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D21').Value = '224.59'
$ws.Range('E21').Value = '  -1.08%  '
$ws.Range('B22').Value = 'Avalanche'
$ws.Range('C22').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D22').Value = '12.33'
$ws.Range('E22').Value = '  -0.86%  '
$ws.Range('B23').Value = 'Dai'
$ws.Range('C23').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D23').Value = '1.000'
$ws.Range('E23').Value = '  -0.03%  '
$ws.Range('B24').Value = 'Chainlink'
$ws.Range('C24').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D24').Value = '7.187'
$ws.Range('E24').Value = '  +0.94%  '
$ws.Range('B25').Value = 'BinanceUSD'
$ws.Range('C25').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D25').Value = '1.000'
$ws.Range('E25').Value = '  +0.05%  '
$ws.Range('B26').Value = 'Monero'
$ws.Range('C26').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D26').Value = '159.39'
$ws.Range('E26').Value = '  +0.57%  '
$ws.Range('B27').Value = 'Cosmos'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D27').Value = '8.386'
$ws.Range('E27').Value = '  -2.32%  '
$ws.Range('B28').Value = 'Stellar'
$ws.Range('C28').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D28').Value = '0.1355'
$ws.Range('E28').Value = '  -2.60%  '
$ws.Range('B29').Value = 'EthereumClassic'
$ws.Range('C29').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D29').Value = '17.81'
$ws.Range('E29').Value = '  -0.55%  '
$ws.Range('B30').Value = 'PancakeSwap'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D30').Value = '1.491'
$ws.Range('E30').Value = '  -1.70%  '
$ws.Range('B31').Value = 'InternetComputer(DFINITY)'
$ws.Range('C31').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D31').Value = '4.035'
$ws.Range('E31').Value = '  -0.15%  '
$ws.Range('B32').Value = 'Filecoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D32').Value = '4.051'
$ws.Range('E32').Value = '  -1.43%  '
$ws.Range('B33').Value = 'Toncoin'
$ws.Range('C33').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D33').Value = '1.200'
$ws.Range('E33').Value = '  +0.81%  '
$ws.Range('B34').Value = 'Hedera'
$ws.Range('C34').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D34').Value = '0.05200'
$ws.Range('E34').Value = '  -0.95%  '
$ws.Range('B35').Value = 'LidoDAOToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D35').Value = '1.835'
$ws.Range('E35').Value = '  -1.55%  '
$ws.Range('B36').Value = 'ARBITRUM'
$ws.Range('C36').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D36').Value = '1.151'
$ws.Range('E36').Value = '  +0.90%  '
$ws.Range('B37').Value = 'ImmutableX'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D37').Value = '0.7314'
$ws.Range('E37').Value = '  -0.82%  '
$ws.Range('B38').Value = 'HuobiToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D38').Value = '2.608'
$ws.Range('E38').Value = '  -1.69%  '
$ws.Range('B39').Value = 'Maker'
$ws.Range('C39').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D39').Value = '1.282.40'
$ws.Range('E39').Value = '  -1.22%  '
$ws.Range('B40').Value = 'MXToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D40').Value = '2.759'
$ws.Range('E40').Value = '  +1.04%  '
$ws.Range('B41').Value = 'VeChain'
$ws.Range('C41').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D41').Value = '0.01790'
$ws.Range('E41').Value = '  -0.01%  '
$ws.Range('B42').Value = 'FraxShare'
$ws.Range('C42').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D42').Value = '6.368'
$ws.Range('E42').Value = '  +6.89%  '
$ws.Range('B43').Value = 'TrustWalletToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D43').Value = '0.8927'
$ws.Range('E43').Value = '  -2.86%  '
$ws.Range('B44').Value = 'PaxDollar'
$ws.Range('C44').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D44').Value = '1.001'
$ws.Range('E44').Value = '  +0.18%  '
$ws.Range('B45').Value = 'Quant'
$ws.Range('C45').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D45').Value = '101.17'
$ws.Range('E45').Value = '  -0.90%  '
$ws.Range('B46').Value = 'RocketPoolETH'
$ws.Range('C46').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D46').Value = '1.980.28'
$ws.Range('E46').Value = '  -0.35%  '
$ws.Range('B47').Value = 'Mantle'
$ws.Range('C47').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D47').Value = '0.5113'
$ws.Range('E47').Value = '  -0.45%  '
$ws.Range('B48').Value = 'Aave'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D48').Value = '63.51'
$ws.Range('E48').Value = '  +0.61%  '
$ws.Range('B49').Value = 'BabyDogeCoin'
$ws.Range('C49').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D49').Value = '0.00000000119'
$ws.Range('E49').Value = '  -0.44%  '
$ws.Range('B50').Value = 'TheSandbox'
$ws.Range('C50').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D50').Value = '0.3977'
$ws.Range('E50').Value = '  -0.56%  '
$ws.Range('B51').Value = 'XinFinNetwork'
$ws.Range('C51').Value = 'https://coinranking.com/coin/77jGXSqWJ1ofG+xinfinnetwork-xdc'
$ws.Range('D51').Value = '0.07231'
$ws.Range('E51').Value = '  -15.87%  '
